# Applies the "Updated symbol list" edit: refreshed Price(D) quotes for the
# existing rows, and a new coin ("One") spliced into the table at row 10,
# shifting WazirX..CoinExToken down one row each (rows 10-18).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Coin / Link / Volume(1h) text columns (rows 10-18 shift down by one) ---
$ws.Range('B10').Value = 'One'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('E10').Value = '9OneONEWorstin24h'
$ws.Range('B11').Value = 'WazirX'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('E11').Value = '10WazirXWRX'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('E12').Value = '11MandalaExchangeTokenMDX'
$ws.Range('B13').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C13').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('E13').Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('B14').Value = 'BitrueCoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('E14').Value = '13BitrueCoinBTR'
$ws.Range('B15').Value = 'BitMartToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('E15').Value = '14BitMartTokenBMX'
$ws.Range('B16').Value = 'MCDex'
$ws.Range('C16').Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range('E16').Value = '15MCDexMCB'
$ws.Range('B17').Value = 'BitForexToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('E17').Value = '16BitForexTokenBF'
$ws.Range('B18').Value = 'CoinExToken'
$ws.Range('C18').Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range('E18').Value = '17CoinExTokenCET'

# --- Price column (D): values are stored as text in this sheet, so force a
#     text number format while writing, then restore each cell's original
#     style so no formatting changes leak into the saved file. ---
$priceCells = @('D2', 'D4', 'D5', 'D7', 'D8', 'D9', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D24', 'D26', 'D40', 'D41', 'D42', 'D44', 'D46', 'D48', 'D49', 'D50')
$priceValues = @('262.83', '6.203', '0.06125', '6.707', '1.363', '0.7985', '0.0006201', '0.1566', '0.08115', '0.03317', '0.03145', '0.09274', '3.935', '0.001693', '0.04816', '0.006178', '0.001099', '0.003375', '0.0001500', '3.694', '2.285', '0.1252', '0.04599', '0.007191', '0.003900', '0.01084', '0.00005978', '0.7001', '0.04939', '0.00002100')
$origStyles = @{}
for ($i = 0; $i -lt $priceCells.Length; $i++) {
    $cell = $ws.Range($priceCells[$i])
    $origStyles[$priceCells[$i]] = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $priceValues[$i]
}
foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = $origStyles[$ref]
}
